$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "3:34PM 12-29-2017"
$ws.Range("B14").Value = "5:04PM 12-29-2017"
$ws.Range("C14").Value = 90

$ws.Range("A15").Value = "6:01PM 12-30-2017"

$ws.Range("B15").Select()
